# Apply the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. '0.558') need an
# explicit text-prefix so Excel keeps storing them as text, matching the
# original inlineStr cells, instead of silently converting to a number.

$ws.Range("D2").Value = "34.535.92"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "1.821.60"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "'0.558"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'35.00"
$ws.Range("E8").Value = "  +8.40%  "
$ws.Range("D9").Value = "'0.299"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").Value = "'0.0695"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "2.083.15"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").Value = "'11.40"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").Value = "1.818.99"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "'0.647"
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").Value = "34.517.86"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").Value = "'4.33"
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").Value = "'69.26"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").Value = "0.0₃0802"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "'247.08"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "'11.57"
$ws.Range("E21").Value = "  +5.28%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'4.18"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").Value = "'169.45"
$ws.Range("E24").Value = "  +4.85%  "
$ws.Range("D25").Value = "'2.08"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("D26").Value = "'7.38"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("D27").Value = "'16.80"
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  +7.13%  "
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").Value = "'1.86"
$ws.Range("E34").Value = "  +2.54%  "
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("D36").Value = "1.419.97"
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("D37").Value = "'0.681"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("D38").Value = "'87.05"
$ws.Range("E38").Value = "  +6.06%  "
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("E41").Value = "  +4.98%  "
$ws.Range("D42").Value = "'0.961"
$ws.Range("E42").Value = "  +4.15%  "
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("D44").Value = "'13.99"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("E46").Value = "  +3.17%  "
$ws.Range("D47").Value = "'6.08"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "1.983.89"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").Value = "'106.07"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0130"
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.20%  "
